# Add a new bulleted/numbered list item "Tests added for Field.cs" right
# after the current last paragraph of the document, matching the style
# (pStyle "a3" + numPr ilvl=0/numId=1) used by the other list items, and
# move the hidden "_GoBack" bookmark (which marks the end of the document /
# last edit point) so that it trails the newly typed run, exactly as real
# Word does when you place the cursor at the end of the document and type
# a new line.

$d = $word.ActiveDocument

# 1. Drop the existing _GoBack bookmark - it currently sits at the end of
#    the last paragraph and we are about to grow the document past it.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Split off a new paragraph after the current last paragraph. Word
#    carries the paragraph style/numbering of the paragraph the break was
#    inserted from, which is exactly the "a3" list style we want.
$lastPara = $d.Paragraphs.Last
$splitPoint = $lastPara.Range
$splitPoint.Collapse(0)
$splitPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.InsertAfter("Tests added for Field.cs")

# 3. This runtime's Bookmarks.Add mis-handles a zero-length range placed
#    exactly at "end of the last paragraph of the document" (and, more
#    generally, right next to an empty paragraph) - it ends up not adding
#    the bookmark where asked. Work around it by temporarily appending a
#    throw-away paragraph with real text, anchoring the bookmark at the
#    boundary right after our new text (which is a perfectly safe spot
#    once it is no longer document-final / paragraph-empty), and then
#    deleting the scratch text back out again. The bookmark keeps tracking
#    that boundary through the clean-up deletes, ending up precisely after
#    "Tests added for Field.cs" - matching real Word's "_GoBack follows the
#    last edit" behaviour.
$afterNewPara = $d.Paragraphs.Last
$scratchStart = $afterNewPara.Range
$scratchStart.Collapse(0)
$scratchStart.InsertParagraphAfter()

$scratchPara = $d.Paragraphs.Last
$scratchPara.Range.InsertAfter("ZZZZ")

$scratchPara = $d.Paragraphs.Last
$anchorPos = $scratchPara.Range.Start

$bookmarkRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove all but the first scratch character for now (deleting a range
# that starts exactly on the bookmark's position would delete the
# bookmark too).
$tailRange = $d.Range($anchorPos + 1, $scratchPara.Range.End)
$tailRange.Delete()

# Merge the one remaining scratch paragraph back into "Tests added for
# Field.cs" by removing the paragraph mark that precedes it (this delete
# sits strictly before the bookmark's position, so the bookmark survives).
$markRange = $d.Range($anchorPos - 1, $anchorPos)
$markRange.Delete()

# Finally, remove the one leftover scratch character right after the
# bookmark.
$leftoverRange = $d.Range($anchorPos - 1, $anchorPos)
$leftoverRange.Delete()
